$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New activity-log entry for 3/14/2019 (row 3), mirroring the existing
# row-2 date formatting.
$ws.Range("A3").Value = 43538
$ws.Range("A3").NumberFormat = "m/d/yy"
$ws.Range("B3").Value = "Downloaded the requried software and configured eclipse to run this project."

# Column B grows to fit the longer activity text (was 41.21875).
$ws.Columns("B:B").ColumnWidth = 62

# Cursor ends up on D4 after the edit, same as the saved workbook.
$ws.Range("D4").Select()
